$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.481.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.404.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.16%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.57"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("E9").Value = "  +6.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.83"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.361"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.58"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.835.70"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.460.71"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000139"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.405.43"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.63%  "
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("E21").Value = "  +5.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.169"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.79"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0770"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.23"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.74"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.29"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.31"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.421"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "304.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.96"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0960"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.410"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.04"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.569"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0225"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.96%  "
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.70%  "
